$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 updates
$ws.Range("E2").Value = [double]"24.07000000000032"
$ws.Range("H2").Value = [double]"5.415722071342227e-16"
$ws.Range("K2").Value = [double]"40.74124452340548"
$ws.Range("L2").Value = "[31.0886152160626, 50.393873830748355]"
$ws.Range("M2").Value = [double]"1.77635683940025e-15"
$ws.Range("N2").Value = [double]"3.552713678800501e-15"
$ws.Range("O2").Value = [double]"1.842816111114733"
$ws.Range("P2").Value = "[1.591237119836272, 2.094395102393195]"
$ws.Range("S2").Value = [double]"63.50221131342981"
$ws.Range("T2").Value = "[58.058628386021496, 68.94579424083813]"
$ws.Range("W2").Value = [double]"17.01043043043066"
$ws.Range("X2").Value = [double]"16.04666666666688"
$ws.Range("Y2").Value = [double]"17.97419419419444"

# Row 3 updates
$ws.Range("E3").Value = [double]"23.84000000000029"
$ws.Range("G3").Value = [double]"1.110223024625157e-15"
$ws.Range("H3").Value = [double]"4.418798107960822e-15"
$ws.Range("K3").Value = [double]"44.88261407523098"
$ws.Range("L3").Value = "[33.77427609838365, 55.9909520520783]"
$ws.Range("M3").Value = [double]"7.127631818093505e-14"
$ws.Range("N3").Value = [double]"7.127631818093505e-14"
$ws.Range("O3").Value = [double]"1.46544762419704"
$ws.Range("P3").Value = "[1.1887107337907326, 1.7421845146033474]"
$ws.Range("S3").Value = [double]"59.13842037402922"
$ws.Range("T3").Value = "[52.16947049878978, 66.10737024926866]"
$ws.Range("W3").Value = [double]"18.27971971971994"
$ws.Range("X3").Value = [double]"17.22970970970992"
$ws.Range("Y3").Value = [double]"19.32972972972996"
